$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Replace-InCell($cell, $find, $replace) {
    $cell.Range.Find.Execute($find, $true, $true, $false, $false, $false, $true, 0, $false, $replace, 1)
}

$cell1 = $t.Cell(1,1)
$cell2 = $t.Cell(1,2)

# --- Cell 1 (BLACKNAVY) ---
Replace-InCell $cell1 "K4" "L11"                      # NO
Replace-InCell $cell1 "SURIYANSYAH" "NARTO"           # NAMA
Replace-InCell $cell1 "43" "42"                       # SEPATU
Replace-InCell $cell1 "DP I NAUTIKA" "DP I TEKNIKA"   # KELAS
Replace-InCell $cell1 "52" "46"                       # UBN_1
Replace-InCell $cell1 "61" "60"                       # UBN_2
Replace-InCell $cell1 "22" "20"                       # UBN_3
Replace-InCell $cell1 "16" "15"                       # UBN_4
Replace-InCell $cell1 "130" "110"                     # UBN_5
Replace-InCell $cell1 "125" "103"                     # UBN_6
Replace-InCell $cell1 "132" "112"                     # UBN_7
Replace-InCell $cell1 "77" "73"                       # UBN_8
Replace-InCell $cell1 "45" "43"                       # UBN_9

# --- Cell 2 (BAJU PUTIH) ---
Replace-InCell $cell2 "K4" "L11"                      # NO
Replace-InCell $cell2 "SURIYANSYAH" "NARTO"           # NAMA
Replace-InCell $cell2 "43" "42"                       # SEPATU
Replace-InCell $cell2 "DP I NAUTIKA" "DP I TEKNIKA"   # KELAS
Replace-InCell $cell2 "49" "46"                       # UH_1
# UH_2 (61) unchanged
Replace-InCell $cell2 "21" "19"                       # UH_3
# UH_4 (13) unchanged
Replace-InCell $cell2 "33" "30"                       # UH_5
Replace-InCell $cell2 "32" "28"                       # UH_6
Replace-InCell $cell2 "33" "30"                       # UH_7 (second occurrence)
Replace-InCell $cell2 "76" "72"                       # UH_8
Replace-InCell $cell2 "45" "43"                       # UH_9
